# Applies the "Update countries & provincias Spain" edit to the Pais sheet:
#  1) Shared-string list reshuffles (Arabia Saudita / Serbia / San Vicente y las
#     Granadinas each move a few slots earlier) which, because the worksheet
#     cells keep referencing the same positional index, changes the country
#     name shown on a handful of rows even though no other data about those
#     rows changed.
#  2) The "Datos actualizados" timestamp advances from 14:50 to 15:20.
#  3) Updated case/recovered/death counts for several countries (new rows get
#     fresh figures; some rows below them simply shift down).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Column A: country-name / header text updates ---
$ws.Cells.Item(1, 1).Value = "Datos actualizados a 3 de Abril de 2020 a las 15:20"
$ws.Cells.Item(39, 1).Value = "Arabia Saudita"
$ws.Cells.Item(40, 1).Value = "Indonesia"
$ws.Cells.Item(41, 1).Value = "Tailandia"
$ws.Cells.Item(45, 1).Value = "Serbia"
$ws.Cells.Item(46, 1).Value = "Panama"
$ws.Cells.Item(47, 1).Value = "Sudafrica"
$ws.Cells.Item(48, 1).Value = "Peru"
$ws.Cells.Item(49, 1).Value = "Republica Dominicana"
$ws.Cells.Item(50, 1).Value = "Islandia"
$ws.Cells.Item(51, 1).Value = "Argentina"
$ws.Cells.Item(206, 1).Value = "San Vicente y las Granadinas"
$ws.Cells.Item(207, 1).Value = "Sierra Leona"
$ws.Cells.Item(208, 1).Value = "Bonaire, San Eustaquio y Saba"

# --- Columns B-H: updated numeric data ---
$ws.Cells.Item(11, 2).Value = 38168
$ws.Cells.Item(11, 3).Value = 4450
$ws.Cells.Item(11, 4).Value = 135
$ws.Cells.Item(11, 5).Value = 34428
$ws.Cells.Item(11, 6).Value = 163
$ws.Cells.Item(11, 7).Value = 684
$ws.Cells.Item(11, 8).Value = 3605
$ws.Cells.Item(16, 2).Value = 11383
$ws.Cells.Item(16, 3).Value = 254
$ws.Cells.Item(16, 4).Value = 2022
$ws.Cells.Item(16, 5).Value = 9193
$ws.Cells.Item(16, 6).Value = 245
$ws.Cells.Item(20, 5).Value = 7620
$ws.Cells.Item(20, 6).Value = 296
$ws.Cells.Item(20, 7).Value = 5
$ws.Cells.Item(20, 8).Value = 329
$ws.Cells.Item(39, 2).Value = 2039
$ws.Cells.Item(39, 3).Value = 154
$ws.Cells.Item(39, 4).Value = 351
$ws.Cells.Item(39, 5).Value = 1663
$ws.Cells.Item(39, 6).Value = 41
$ws.Cells.Item(39, 7).Value = 4
$ws.Cells.Item(39, 8).Value = 25
$ws.Cells.Item(40, 2).Value = 1986
$ws.Cells.Item(40, 3).Value = 196
$ws.Cells.Item(40, 4).Value = 134
$ws.Cells.Item(40, 5).Value = 1671
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 11
$ws.Cells.Item(40, 8).Value = 181
$ws.Cells.Item(41, 2).Value = 1978
$ws.Cells.Item(41, 3).Value = 103
$ws.Cells.Item(41, 4).Value = 581
$ws.Cells.Item(41, 5).Value = 1378
$ws.Cells.Item(41, 6).Value = 23
$ws.Cells.Item(41, 7).Value = 4
$ws.Cells.Item(41, 8).Value = 19
$ws.Cells.Item(45, 2).Value = 1476
$ws.Cells.Item(45, 3).Value = 305
$ws.Cells.Item(45, 4).Value = 42
$ws.Cells.Item(45, 5).Value = 1395
$ws.Cells.Item(45, 6).Value = 81
$ws.Cells.Item(45, 7).Value = 8
$ws.Cells.Item(45, 8).Value = 39
$ws.Cells.Item(46, 2).Value = 1475
$ws.Cells.Item(46, 3).Value = 0
$ws.Cells.Item(46, 4).Value = 9
$ws.Cells.Item(46, 5).Value = 1429
$ws.Cells.Item(46, 6).Value = 50
$ws.Cells.Item(46, 7).Value = 0
$ws.Cells.Item(46, 8).Value = 37
$ws.Cells.Item(47, 2).Value = 1462
$ws.Cells.Item(47, 3).Value = 0
$ws.Cells.Item(47, 4).Value = 95
$ws.Cells.Item(47, 5).Value = 1362
$ws.Cells.Item(47, 6).Value = 7
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 8).Value = 5
$ws.Cells.Item(48, 2).Value = 1414
$ws.Cells.Item(48, 3).Value = 0
$ws.Cells.Item(48, 4).Value = 537
$ws.Cells.Item(48, 5).Value = 822
$ws.Cells.Item(48, 6).Value = 51
$ws.Cells.Item(48, 7).Value = 0
$ws.Cells.Item(48, 8).Value = 55
$ws.Cells.Item(49, 2).Value = 1380
$ws.Cells.Item(49, 3).Value = 0
$ws.Cells.Item(49, 4).Value = 16
$ws.Cells.Item(49, 5).Value = 1304
$ws.Cells.Item(49, 6).Value = 147
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = 60
$ws.Cells.Item(50, 2).Value = 1364
$ws.Cells.Item(50, 3).Value = 45
$ws.Cells.Item(50, 4).Value = 309
$ws.Cells.Item(50, 5).Value = 1051
$ws.Cells.Item(50, 6).Value = 12
$ws.Cells.Item(50, 7).Value = 0
$ws.Cells.Item(50, 8).Value = 4
$ws.Cells.Item(51, 2).Value = 1265
$ws.Cells.Item(51, 3).Value = 0
$ws.Cells.Item(51, 4).Value = 256
$ws.Cells.Item(51, 5).Value = 972
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 1
$ws.Cells.Item(51, 8).Value = 37
$ws.Cells.Item(71, 2).Value = 568
$ws.Cells.Item(71, 3).Value = 35
$ws.Cells.Item(71, 4).Value = 27
$ws.Cells.Item(71, 5).Value = 524
$ws.Cells.Item(71, 6).Value = 4
$ws.Cells.Item(73, 4).Value = 26
$ws.Cells.Item(73, 5).Value = 473
$ws.Cells.Item(73, 6).Value = 65
$ws.Cells.Item(206, 2).Value = 3
$ws.Cells.Item(206, 3).Value = 1
$ws.Cells.Item(206, 4).Value = 1
$ws.Cells.Item(206, 5).Value = 2
$ws.Cells.Item(208, 4).Value = 0
$ws.Cells.Item(208, 5).Value = 2
$ws.Cells.Item(208, 6).Value = 0
